# edit.ps1 -- apply the CV.docx changes described by the commit diff:
#   1) Add a new "Compact" bullet paragraph (AWS summary bullet) right after
#      the "Providing in-house support and training to other team members"
#      bullet in the CDL / DevOps Engineer section.
#   2) Remove the "References" Heading2 paragraph and its
#      "References available upon request" body paragraph (the bookmark
#      that wrapped them, id 35 "references", is collapsed along with it;
#      the outer "experience-summary" bookmark's id shifts down from 36 to
#      35 to take its place once the runtime re-serialises the bookmark
#      table).

$d = $word.ActiveDocument

function Get-ParagraphText($para) {
    # Paragraph.Range.Text includes the trailing paragraph mark (CR / 0x0D,
    # sometimes a cell mark 0x07) -- strip that off before comparing.
    return $para.Range.Text.TrimEnd([char]13, [char]7)
}

function Find-ParagraphByText($doc, [string]$text) {
    $paras = $doc.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        if ((Get-ParagraphText $p) -eq $text) {
            return $p
        }
    }
    return $null
}

function Find-ParagraphIndexByText($doc, [string]$text) {
    $paras = $doc.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        if ((Get-ParagraphText $p) -eq $text) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------
# 1) Insert the new AWS bullet paragraph after the "Providing in-house
#    support..." bullet that belongs to the CDL / DevOps Engineer role
#    (the first occurrence of that sentence -- Laterooms has a similar but
#    distinct sentence later in the document).
# ---------------------------------------------------------------------
$anchor = Find-ParagraphByText $d "Providing in-house support and training to other team members"
if ($anchor -eq $null) {
    throw "Could not find the 'Providing in-house support...' anchor paragraph"
}

# InsertParagraphAfter on the anchor's range creates a new, empty sibling
# paragraph that inherits the anchor's paragraph formatting (style Compact +
# the numPr list numbering), exactly like pressing Enter at the end of the
# bullet in Word.
$anchor.Range.InsertParagraphAfter()

$anchorIndex = Find-ParagraphIndexByText $d "Providing in-house support and training to other team members"
$newPara = $d.Paragraphs.Item($anchorIndex + 1)
$newPara.Range.Text = "AWS use primaraly includes ECS (both Fargate and EC2 backed), S3, EC2, Route53, CodeDeploy, ALB/ELB and to a lesser extent Cloudfront, GlobalAccelerator, Labda. All primaly provisioned via Terrafrom but with knowlage of manual deployement/corrections of services if needed."

# ---------------------------------------------------------------------
# 2) Remove the "References" section: the Heading2 "References" paragraph
#    and the "References available upon request" FirstParagraph that
#    follows it.
# ---------------------------------------------------------------------
$headingIndex = Find-ParagraphIndexByText $d "References"
if ($headingIndex -eq -1) {
    throw "Could not find the 'References' heading paragraph"
}
$bodyIndex = $headingIndex + 1
$headingPara = $d.Paragraphs.Item($headingIndex)
$bodyPara = $d.Paragraphs.Item($bodyIndex)

$deleteStart = $headingPara.Range.Start
$deleteEnd = $bodyPara.Range.End
$d.Range($deleteStart, $deleteEnd).Delete()

"done"
